$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 in the sheet contains a duplicate department entry
# (it repeats the content from row 17). Delete the entire row,
# which shifts all subsequent rows up by one.
$ws.Rows.Item(19).Delete()

# Column A holds a plain sequential index (row number - 2). Restore the
# sequence for the rows that shifted up after the delete.
for ($r = 19; $r -le 35; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
